# Add a new worksheet "Sheet2" after the existing "Sheet1" and populate it
# with policy/quote test data, mirroring the commit
# "data writing to excel code is added".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 so it becomes the 2nd (and active) tab.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# Data rows
# Numbers-as-text cells (e.g. "11538380") must be stored as shared strings
# (t="s"), not numeric values, to match the source data - use a leading
# apostrophe to force text entry, then drop the resulting text-number style
# so the cell is left with no explicit style (matching the target file).
# Write order matches the original authoring sequence (data rows first,
# then the header row, then the final row) so new shared-string entries
# land at the same indexes as the source file.
$ws2.Range("A2").Value = "TN2485269"
$ws2.Range("B2").Value = "'11538380"
$ws2.Range("B2").ClearFormats()

$ws2.Range("A3").Value = "TB2485272"
$ws2.Range("B3").Value = "'11538402"
$ws2.Range("B3").ClearFormats()

# Header row
$ws2.Range("A1").Value = "Policy Number"
$ws2.Range("B1").Value = "Quote Number"
$ws2.Range("C1").Value = "Type Policy"
$ws2.Range("D1").Value = "Type"
$ws2.Range("E1").Value = "Testcaseid"

$ws2.Range("A4").Value = "TB2485273"
$ws2.Range("B4").Value = "'11538440"
$ws2.Range("B4").ClearFormats()
$ws2.Range("C4").Value = "Bond - No Credit"
$ws2.Range("D4").Value = "TC001"

# Column widths for the new sheet (closest achievable values given the
# host's column-width quantization; target stored widths are 16.140625,
# 15.0, 16.7109375 and 15.42578125 respectively).
$ws2.Columns.Item(1).ColumnWidth = 15.333333333333334
$ws2.Columns.Item(2).ColumnWidth = 14.166666666666666
$ws2.Columns.Item(3).ColumnWidth = 15.833333333333334
$ws2.Columns.Item(5).ColumnWidth = 14.666666666666666

# Selection on the new sheet
$ws2.Range("I8").Select() | Out-Null

$ws2.Activate() | Out-Null
